# B6-PowerPoint.pptx edit script
# 1) Re-style the three data tables (slides 14, 15, 16) with the new
#    built-in table style GUID.
# 2) Swap the deck's theme palette from "Integral / Red Violet" to the
#    stock "Office" colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink),
#    mirroring the theme1.xml <-> theme2.xml content swap.

$p = $ppt.ActivePresentation

# --- 1. Table style updates -------------------------------------------------
$newStyleId = "{4E98B519-F85A-4075-9106-43047619D924}"

foreach ($slideIdx in @(14, 15, 16)) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour scheme updates -----------------------------------------
# Target "Office" theme colours (RGB hex -> BGR-encoded long used by the
# PowerPoint object model's .RGB property).
#   index : role      : hex    : BGR-encoded value
#     1   : dk1        000000   0
#     2   : lt1         FFFFFF   16777215
#     3   : dk2         44546A   6968388
#     4   : lt2         E7E6E6   15132391
#     5   : accent1     5B9BD5   13998939
#     6   : accent2     ED7D31   3243501
#     7   : accent3     A5A5A5   10855845
#     8   : accent4     FFC000   49407
#     9   : accent5     4472C4   12874308
#    10   : accent6     70AD47   4697456
#    11   : hlink       0563C1   12673797
#    12   : folHlink    954F72   7491477
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$design = $p.Designs.Item(1)
$themeColorScheme = $design.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColorScheme.Item($i).RGB = $officeColors[$i - 1]
}
